# Commit: "Commit to add the file: PdpNavigation.java"
# The second worksheet (ProductCatalogueData) becomes a login-only data sheet
# (LoginOnly) used by the new PdpNavigation test: one header row plus three
# login-credential rows. The previously 4th data row (old row 5) is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# 1. Rename the sheet.
$ws.Name = "LoginOnly"

# 2. Remove the last data row (old row 5) -- select it first so the
#    resulting selection matches a genuine "select row, delete row" edit.
$ws.Rows.Item(5).Select() | Out-Null
$ws.Rows.Item(5).Delete() | Out-Null

# 3. Refresh the hyperlink-bearing cells with the new login/URL data.
#    Clearing first makes sure stale rId-based hyperlinks don't linger.
$ws.Hyperlinks.Delete()

$ws.Range("C2").Value = "https://www.everlast.com/"
$ws.Range("A2").Value = "blubd.softtech@gmail.com"
$ws.Range("B2").Value = "demotest"

$ws.Range("A3").Value = "wrongId@wrong.com"
$ws.Range("B3").Value = "demotest"
$ws.Range("C3").Value = "https://www.everlast.com/"

$ws.Range("A4").Value = "blubd.softtech@gmail.com"
$ws.Range("B4").Value = "everlast#123"
$ws.Range("C4").Value = "https://www.everlast.com/"

# 4. Re-create the hyperlinks: the loginUrl column (C) plus the two user
#    name cells that now carry e-mail addresses (A2, A3).
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.everlast.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.everlast.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.everlast.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:blubd.softtech@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:wrongId@wrong.com") | Out-Null

# 5. Re-size the columns to fit the new, wider content (approximates the
#    bestFit autofit widths Excel would have computed for the new values).
$ws.Columns.Item(1).ColumnWidth = 22.585
$ws.Columns.Item(2).ColumnWidth = 8.585
$ws.Columns.Item(3).ColumnWidth = 21.752
